# Update the "想去人数" (want-to-go count, column F) figures across all
# four worksheets (展览, 演出, 本地生活, 全部类型) to match the latest
# scrape of the source data (gh-pages output generated at 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws = $wb.Worksheets.Item(1)
$ws.Range("F5").Value = 7423
$ws.Range("F7").Value = 4775
$ws.Range("F8").Value = 7024
$ws.Range("F13").Value = 174
$ws.Range("F15").Value = 1163
$ws.Range("F17").Value = 160
$ws.Range("F19").Value = 224
$ws.Range("F21").Value = 1157
$ws.Range("F22").Value = 952
$ws.Range("F30").Value = 175
$ws.Range("F32").Value = 38
$ws.Range("F33").Value = 92
$ws.Range("F37").Value = 69
$ws.Range("F39").Value = 371
$ws.Range("F42").Value = 139
$ws.Range("F43").Value = 22

# Sheet 2: 演出
$ws = $wb.Worksheets.Item(2)
$ws.Range("F12").Value = 28
$ws.Range("F17").Value = 556
$ws.Range("F21").Value = 203
$ws.Range("F26").Value = 635
$ws.Range("F28").Value = 26
$ws.Range("F32").Value = 987
$ws.Range("F33").Value = 609
$ws.Range("F35").Value = 1

# Sheet 3: 本地生活
$ws = $wb.Worksheets.Item(3)
$ws.Range("F4").Value = 729
$ws.Range("F5").Value = 854
$ws.Range("F6").Value = 671
$ws.Range("F8").Value = 1583
$ws.Range("F9").Value = 2479

# Sheet 4: 全部类型
$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 729
$ws.Range("F4").Value = 854
$ws.Range("F7").Value = 671
$ws.Range("F8").Value = 671
$ws.Range("F9").Value = 7423
$ws.Range("F11").Value = 4775
$ws.Range("F13").Value = 7024
$ws.Range("F17").Value = 174
$ws.Range("F18").Value = 1583
$ws.Range("F19").Value = 2479
$ws.Range("F20").Value = 203
$ws.Range("F22").Value = 1163
$ws.Range("F23").Value = 160
$ws.Range("F26").Value = 1157
$ws.Range("F27").Value = 635
$ws.Range("F28").Value = 952
$ws.Range("F32").Value = 175
$ws.Range("F33").Value = 26
$ws.Range("F35").Value = 38
$ws.Range("F36").Value = 92
$ws.Range("F37").Value = 987
$ws.Range("F39").Value = 609
$ws.Range("F40").Value = 69
$ws.Range("F43").Value = 371
$ws.Range("F48").Value = 139
$ws.Range("F50").Value = 22

$wb.Save()
